$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Add the three new rows of data (rows 94-96) ----

# Row 94: 217. Contains Duplicate
$ws.Range("A94").Value = "217. Contains Duplicate"
$ws.Range("B94").Value = "Easy"
$ws.Range("C94").Value = "https://leetcode.com/problems/contains-duplicate/"
$ws.Range("D94").Value = 44591
$ws.Range("E94").Value = "哈希"
$ws.Range("F94").Value = "简单用set查重即可"
$ws.Range("G93").Copy($ws.Range("G94"))

# Row 95: 697. Degree of an Array
$ws.Range("A95").Value = "697. Degree of an Array"
$ws.Range("B95").Value = "Easy"
$ws.Range("C95").Value = "https://leetcode.com/problems/degree-of-an-array/"
$ws.Range("D95").Value = 44591
$ws.Range("E95").Value = "哈希"
$ws.Range("F95").Value = "hash需要记录两个数据，一个是出现的次数，一个是序列的长度"
$ws.Range("G93").Copy($ws.Range("G95"))

# Row 96: 594. Longest Harmonious Subsequence
$ws.Range("A96").Value = "594. Longest Harmonious Subsequence"
$ws.Range("B96").Value = "Easy"
$ws.Range("C96").Value = "https://leetcode.com/problems/longest-harmonious-subsequence/"
$ws.Range("D96").Value = 44591
$ws.Range("E96").Value = "哈希"
$ws.Range("F96").Value = "关键就是最后的序列里面只有相邻的两个值"
$ws.Range("G93").Copy($ws.Range("G96"))
$ws.Range("H65").Copy($ws.Range("H96"))

# ---- Hyperlinks for the new problem links ----
$ws.Hyperlinks.Add($ws.Range("C94"), "https://leetcode.com/problems/contains-duplicate/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C95"), "https://leetcode.com/problems/degree-of-an-array/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C96"), "https://leetcode.com/problems/longest-harmonious-subsequence/") | Out-Null

# ---- Fix up formatting to match the rest of the table ----
# Rows 94 & 95 use the same banding/format as rows 91/92
$ws.Range("A91:F92").Copy()
$ws.Range("A94:F95").PasteSpecial(-4122)

# Row 96 uses the alternate banding/format (same as row 65 which has the
# same pattern including the H "error mark" column)
$ws.Range("A65:H65").Copy()
$ws.Range("A96:H96").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Match the row height used throughout the rest of the table
$ws.Rows(94).RowHeight = 28
$ws.Rows(95).RowHeight = 28
$ws.Rows(96).RowHeight = 28

# ---- View tweaks: keep selection in sync with the newly added rows ----
$ws.Range("F107").Select() | Out-Null
